# Update crypto price/volume data per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need NumberFormat forced to
# Text ("@") first, otherwise the COM Value setter auto-converts the string to
# a numeric value (losing the trailing zero / exact text representation).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"

$ws.Range("D2").Value = "38.819.70"
$ws.Range("E2").Value = "  +2.73%  "
$ws.Range("D3").Value = "2.092.39"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "228.32"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("D7").Value = "60.84"
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").Value = "0.0837"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "2.403.03"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("E13").Value = "  +3.57%  "
$ws.Range("D14").Value = "22.03"
$ws.Range("E14").Value = "  +4.77%  "
$ws.Range("D15").Value = "0.798"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "2.106.48"
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("D18").Value = "38.750.63"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("D19").Value = "71.82"
$ws.Range("E19").Value = "  +3.28%  "
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").Value = "170.70"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "9.42"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  +6.33%  "
$ws.Range("D29").Value = "1.46"
$ws.Range("E29").Value = "  +13.03%  "
$ws.Range("D30").Value = "19.17"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").Value = "4.48"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  +4.71%  "
$ws.Range("D35").Value = "0.0613"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "2.40"
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "6.44"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("D38").Value = "3.52"
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "18.29"
$ws.Range("E40").Value = "  +1.56%  "
$ws.Range("D41").Value = "101.25"
$ws.Range("E41").Value = "  +3.72%  "
$ws.Range("D42").Value = "1.538.67"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("E43").Value = "  +3.44%  "
$ws.Range("D44").Value = "0.0926"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("E46").Value = "  +7.84%  "
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("E51").Value = "  +2.44%  "

# Restore the default (General) style on those cells so no residual text
# formatting is left behind - matches original workbook styling.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
